$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update price (D) and 1h volume/change (E) columns for rows whose values changed
$ws.Range("D2").Value = "22.455.79"
$ws.Range("E2").Value = "  +0.23%  "
$ws.Range("D3").Value = "1.572.14"
$ws.Range("E3").Value = "  +0.05%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("E5").Value = "  +0.04%  "
$ws.Range("D6").Value = "'291.72"
$ws.Range("E6").Value = "  +0.00%  "
$ws.Range("D7").Value = "'0.3722"
$ws.Range("E7").Value = "  -1.03%  "
$ws.Range("D8").Value = "'49.78"
$ws.Range("E8").Value = "  -0.24%  "
$ws.Range("D9").Value = "'0.3407"
$ws.Range("E9").Value = "  -0.42%  "
$ws.Range("E10").Value = "  -0.01%  "
$ws.Range("D11").Value = "'0.07552"
$ws.Range("E11").Value = "  -1.26%  "
$ws.Range("D12").Value = "'1.002"
$ws.Range("E12").Value = "  +0.01%  "
$ws.Range("D13").Value = "'21.26"
$ws.Range("E13").Value = "  +0.24%  "
$ws.Range("D14").Value = "'6.037"
$ws.Range("E14").Value = "  +0.53%  "
$ws.Range("D15").Value = "'6.965"
$ws.Range("E15").Value = "  +0.24%  "
$ws.Range("D16").Value = "1.571.04"
$ws.Range("E16").Value = "  +0.12%  "
$ws.Range("D17").Value = "'0.00001124"
$ws.Range("E17").Value = "  -0.91%  "
$ws.Range("D18").Value = "'91.27"
$ws.Range("E18").Value = "  +1.26%  "
$ws.Range("D19").Value = "'0.06750"
$ws.Range("E19").Value = "  +0.24%  "
$ws.Range("D20").Value = "'1.001"
$ws.Range("E20").Value = "  +0.00%  "
$ws.Range("E21").Value = "  +1.25%  "
$ws.Range("D22").Value = "'16.33"
$ws.Range("E22").Value = "  -2.52%  "
$ws.Range("E23").Value = "  +1.17%  "
$ws.Range("D24").Value = "22.447.44"
$ws.Range("E24").Value = "  +0.24%  "
$ws.Range("D25").Value = "'2.376"
$ws.Range("E25").Value = "  -0.85%  "
$ws.Range("D26").Value = "'2.681"
$ws.Range("E26").Value = "  +1.06%  "
$ws.Range("D27").Value = "'20.03"
$ws.Range("E27").Value = "  -0.68%  "
$ws.Range("D28").Value = "'148.93"
$ws.Range("E28").Value = "  +1.22%  "
$ws.Range("D29").Value = "'5.053"
$ws.Range("E29").Value = "  +0.52%  "
$ws.Range("D30").Value = "'125.75"
$ws.Range("E30").Value = "  -0.75%  "
$ws.Range("D31").Value = "1.748.08"
$ws.Range("E31").Value = "  +0.21%  "
$ws.Range("D32").Value = "'1.081"
$ws.Range("E32").Value = "  +10.07%  "
$ws.Range("D33").Value = "'6.209"
$ws.Range("E33").Value = "  +0.65%  "
$ws.Range("D34").Value = "'2.016"
$ws.Range("E34").Value = "  +0.57%  "
$ws.Range("D35").Value = "'9.835"
$ws.Range("E35").Value = "  -2.63%  "
$ws.Range("D36").Value = "'0.08380"
$ws.Range("E36").Value = "  -1.43%  "
$ws.Range("D37").Value = "'0.02491"
$ws.Range("E37").Value = "  -1.88%  "
$ws.Range("D40").Value = "'0.06539"
$ws.Range("E40").Value = "  -0.26%  "
$ws.Range("D41").Value = "'5.465"
$ws.Range("E41").Value = "  +0.90%  "
$ws.Range("D42").Value = "'11.38"
$ws.Range("E42").Value = "  -0.70%  "
$ws.Range("D43").Value = "'0.6239"
$ws.Range("E43").Value = "  -2.46%  "
$ws.Range("E44").Value = "  +0.00%  "
$ws.Range("D45").Value = "'14.02"
$ws.Range("E45").Value = "  -0.69%  "
$ws.Range("D46").Value = "'3.813"
$ws.Range("E46").Value = "  +0.27%  "
$ws.Range("D47").Value = "'0.5821"
$ws.Range("E47").Value = "  -2.55%  "
$ws.Range("D48").Value = "'130.39"
$ws.Range("E48").Value = "  +4.49%  "
$ws.Range("D49").Value = "'2.076"
$ws.Range("E49").Value = "  -0.91%  "
$ws.Range("D50").Value = "'1.223"
$ws.Range("E50").Value = "  -4.99%  "
$ws.Range("D51").Value = "'0.07323"
$ws.Range("E51").Value = "  -0.07%  "

# Rows 38 and 39 swapped order: Algorand now ranks above TrustWalletToken
$ws.Range("B38").Value = "Algorand"
$ws.Range("C38").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D38").Value = "'0.2305"
$ws.Range("E38").Value = "  -0.32%  "
$ws.Range("B39").Value = "TrustWalletToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D39").Value = "'1.347"
$ws.Range("E39").Value = "  -2.91%  "

Write-Host "Applied cryptos list update"